$d = $word.ActiveDocument
for ($i=1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  $t = $p.Range.Text
  if ($t.Length -gt 60) { $t = $t.Substring(0,60) }
  Write-Host $i ": [" $t "]"
}
